# Append the new Timer-Results row (row 40) with the latest Perfecto run data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

# Columns that hold numeric-looking identifiers (accounts/cases/contacts/people)
# must be forced to Text first so Excel stores them as literal strings
# ("6167", "8759", ...) instead of silently re-typing them as numbers.
$ws.Range($ws.Cells.Item($row, 3), $ws.Cells.Item($row, 6)).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "https://salesforce.perfectomobile.com/services/reports/PRIVATE:Raj_SalesForceCaseScript_16-08-04_13_37_43_99227.xml?operation=download&format=html&user=rajp@perfectomobile.com&password="
$ws.Cells.Item($row, 2).Value = "04/08/2016 09:37:43"
$ws.Cells.Item($row, 3).Value = "6167"
$ws.Cells.Item($row, 4).Value = "8759"
$ws.Cells.Item($row, 5).Value = "6770"
$ws.Cells.Item($row, 6).Value = "3980"
$ws.Cells.Item($row, 7).Value = "https://demo.vod-download-01.perfectomobile.com/demo/647672315467564c494d713734357762553872464b7a326b6638337439776b4c4c6c5a55703352456c31453d/877772dae0f270833cc141538dfe10f7f1b12bd16b44b08802322d2518240e5e.flv"
$ws.Cells.Item($row, 8).Value = "PRIVATE:Raj_SalesForceCaseScript_16-08-04_13_37_43_99227.xml"
